$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42
$ws.Range("B42").Value = 6149872
$ws.Range("E42").Value = "FC Seoul"
$ws.Range("F42").Value = "Suwon FC"
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 2
$ws.Range("I42").Value = "H"
$ws.Range("J42").Value = 1.55
$ws.Range("K42").Value = 3.75
$ws.Range("L42").Value = 5.25
$ws.Range("M42").Value = 1.55
$ws.Range("N42").Value = 4
$ws.Range("O42").Value = 5
$ws.Range("P42").Value = -1
$ws.Range("Q42").Value = 1.925
$ws.Range("R42").Value = 1.925
$ws.Range("S42").Value = 3
$ws.Range("T42").Value = 1.9
$ws.Range("U42").Value = 1.95
$ws.Range("V42").Value = 0.55
$ws.Range("W42").Value = -1
$ws.Range("Y42").Value = 0.925
$ws.Range("AA42").Value = 0.8999999999999999
$ws.Range("AB42").Value = -1

# Row 43
$ws.Range("B43").Value = 6149871
$ws.Range("E43").Value = "Suwon Bluewings"
$ws.Range("F43").Value = "Pohang Steelers"
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 1
$ws.Range("I43").Value = "D"
$ws.Range("J43").Value = 4
$ws.Range("K43").Value = 3.3
$ws.Range("L43").Value = 1.909
$ws.Range("M43").Value = 4.5
$ws.Range("N43").Value = 3.3
$ws.Range("O43").Value = 1.833
$ws.Range("P43").Value = 0.5
$ws.Range("Q43").Value = 2.025
$ws.Range("R43").Value = 1.825
$ws.Range("S43").Value = 2.25
$ws.Range("T43").Value = 1.925
$ws.Range("U43").Value = 1.925
$ws.Range("V43").Value = -1
$ws.Range("W43").Value = 2.3
$ws.Range("Y43").Value = 1.025
$ws.Range("AA43").Value = -0.5
$ws.Range("AB43").Value = 0.4625

# Row 56
$ws.Range("B56").Value = 6149879
$ws.Range("E56").Value = "Gwangju FC"
$ws.Range("F56").Value = "Daejeon Hana Citizen"
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = "H"
$ws.Range("J56").Value = 2.15
$ws.Range("K56").Value = 3.3
$ws.Range("L56").Value = 3.4
$ws.Range("M56").Value = 2.2
$ws.Range("N56").Value = 3
$ws.Range("O56").Value = 3.6
$ws.Range("P56").Value = -0.25
$ws.Range("Q56").Value = 1.9
$ws.Range("R56").Value = 1.95
$ws.Range("T56").Value = 2
$ws.Range("U56").Value = 1.85
$ws.Range("V56").Value = 1.2
$ws.Range("W56").Value = -1
$ws.Range("Y56").Value = 0.8999999999999999
$ws.Range("Z56").Value = -1
$ws.Range("AA56").Value = 1

# Row 57
$ws.Range("B57").Value = 6149880
$ws.Range("E57").Value = "FC Seoul"
$ws.Range("F57").Value = "Pohang Steelers"
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 2
$ws.Range("I57").Value = "D"
$ws.Range("J57").Value = 2.5
$ws.Range("K57").Value = 3.2
$ws.Range("L57").Value = 2.8
$ws.Range("M57").Value = 2.625
$ws.Range("N57").Value = 3.3
$ws.Range("O57").Value = 2.6
$ws.Range("P57").Value = 0
$ws.Range("Q57").Value = 1.95
$ws.Range("R57").Value = 1.9
$ws.Range("T57").Value = 1.825
$ws.Range("U57").Value = 2.025
$ws.Range("V57").Value = -1
$ws.Range("W57").Value = 2.3
$ws.Range("Y57").Value = 0
$ws.Range("Z57").Value = 0
$ws.Range("AA57").Value = 0.825

# Row 68
$ws.Range("B68").Value = 6149886
$ws.Range("E68").Value = "Suwon Bluewings"
$ws.Range("F68").Value = "Jeju United"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = "H"
$ws.Range("J68").Value = 2.9
$ws.Range("K68").Value = 3.25
$ws.Range("L68").Value = 2.2
$ws.Range("M68").Value = 2.75
$ws.Range("N68").Value = 3.2
$ws.Range("O68").Value = 2.375
$ws.Range("Q68").Value = 2.1
$ws.Range("R68").Value = 1.775
$ws.Range("S68").Value = 2.5
$ws.Range("T68").Value = 2.05
$ws.Range("U68").Value = 1.8
$ws.Range("V68").Value = 1.75
$ws.Range("W68").Value = -1
$ws.Range("Y68").Value = 1.1
$ws.Range("Z68").Value = -1
$ws.Range("AA68").Value = -1
$ws.Range("AB68").Value = 0.8

# Row 69
$ws.Range("B69").Value = 6149887
$ws.Range("E69").Value = "Incheon Utd"
$ws.Range("F69").Value = "Gwangju FC"
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 2
$ws.Range("I69").Value = "D"
$ws.Range("J69").Value = 2.375
$ws.Range("K69").Value = 3
$ws.Range("L69").Value = 2.8
$ws.Range("M69").Value = 2.625
$ws.Range("N69").Value = 3
$ws.Range("O69").Value = 2.7
$ws.Range("Q69").Value = 1.9
$ws.Range("R69").Value = 1.95
$ws.Range("S69").Value = 2
$ws.Range("T69").Value = 1.85
$ws.Range("U69").Value = 2
$ws.Range("V69").Value = -1
$ws.Range("W69").Value = 2
$ws.Range("Y69").Value = 0
$ws.Range("Z69").Value = 0
$ws.Range("AA69").Value = 0.8500000000000001
$ws.Range("AB69").Value = -1

# Row 147
$ws.Range("B147").Value = 7715261
$ws.Range("E147").Value = "Pohang Steelers"
$ws.Range("F147").Value = "Daegu FC"
$ws.Range("G147").Value = 3
$ws.Range("H147").Value = 1
$ws.Range("I147").Value = "H"
$ws.Range("J147").Value = 2.05
$ws.Range("K147").Value = 3.3
$ws.Range("L147").Value = 3.6
$ws.Range("M147").Value = 2.75
$ws.Range("N147").Value = 3.1
$ws.Range("O147").Value = 2.7
$ws.Range("P147").Value = 0
$ws.Range("Q147").Value = 1.95
$ws.Range("R147").Value = 1.9
$ws.Range("S147").Value = 2
$ws.Range("T147").Value = 1.85
$ws.Range("U147").Value = 2
$ws.Range("V147").Value = 1.75
$ws.Range("X147").Value = -1
$ws.Range("Y147").Value = 0.95
$ws.Range("Z147").Value = -1
$ws.Range("AA147").Value = 0.8500000000000001

# Row 148
$ws.Range("B148").Value = 7715262
$ws.Range("E148").Value = "Gimcheon Sangmu FC"
$ws.Range("F148").Value = "Ulsan Hyundai"
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 3
$ws.Range("I148").Value = "A"
$ws.Range("J148").Value = 3.3
$ws.Range("K148").Value = 3.5
$ws.Range("L148").Value = 2.05
$ws.Range("M148").Value = 3
$ws.Range("N148").Value = 3.2
$ws.Range("O148").Value = 2.4
$ws.Range("P148").Value = 0.25
$ws.Range("Q148").Value = 1.775
$ws.Range("R148").Value = 2.1
$ws.Range("S148").Value = 2.25
$ws.Range("T148").Value = 2
$ws.Range("U148").Value = 1.85
$ws.Range("V148").Value = -1
$ws.Range("X148").Value = 1.4
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = 1.1
$ws.Range("AA148").Value = 1

# Row 173
$ws.Range("B173").Value = 7716531
$ws.Range("E173").Value = "Daegu FC"
$ws.Range("F173").Value = "FC Seoul"
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = "D"
$ws.Range("J173").Value = 2.45
$ws.Range("K173").Value = 3.2
$ws.Range("L173").Value = 2.7
$ws.Range("M173").Value = 2.5
$ws.Range("N173").Value = 3.25
$ws.Range("O173").Value = 2.8
$ws.Range("P173").Value = 0
$ws.Range("Q173").Value = 1.85
$ws.Range("R173").Value = 2
$ws.Range("T173").Value = 1.9
$ws.Range("U173").Value = 1.95
$ws.Range("W173").Value = 2.25
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = 0
$ws.Range("Z173").Value = 0
$ws.Range("AA173").Value = -1
$ws.Range("AB173").Value = 0.95

# Row 174
$ws.Range("B174").Value = 7715278
$ws.Range("E174").Value = "Jeonbuk Motors"
$ws.Range("F174").Value = "Gangwon FC"
$ws.Range("G174").Value = 2
$ws.Range("H174").Value = 3
$ws.Range("I174").Value = "A"
$ws.Range("J174").Value = 1.8
$ws.Range("K174").Value = 3.5
$ws.Range("L174").Value = 4.75
$ws.Range("M174").Value = 1.909
$ws.Range("N174").Value = 3.4
$ws.Range("O174").Value = 4
$ws.Range("P174").Value = -0.5
$ws.Range("Q174").Value = 1.975
$ws.Range("R174").Value = 1.875
$ws.Range("T174").Value = 1.95
$ws.Range("U174").Value = 1.9
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 3
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = 0.875
$ws.Range("AA174").Value = 0.95
$ws.Range("AB174").Value = -1

# Row 188
$ws.Range("B188").Value = 7716476
$ws.Range("E188").Value = "Gwangju FC"
$ws.Range("F188").Value = "Suwon FC"
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = "A"
$ws.Range("J188").Value = 2.3
$ws.Range("L188").Value = 3.4
$ws.Range("M188").Value = 1.909
$ws.Range("N188").Value = 3.3
$ws.Range("O188").Value = 4.333
$ws.Range("P188").Value = -0.5
$ws.Range("Q188").Value = 1.9
$ws.Range("R188").Value = 1.95
$ws.Range("S188").Value = 2.5
$ws.Range("T188").Value = 1.95
$ws.Range("U188").Value = 1.9
$ws.Range("V188").Value = -1
$ws.Range("X188").Value = 3.333
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = 0.95
$ws.Range("AA188").Value = 0.95
$ws.Range("AB188").Value = -1

# Row 189
$ws.Range("B189").Value = 7716477
$ws.Range("E189").Value = "Gimcheon Sangmu FC"
$ws.Range("F189").Value = "Gangwon FC"
$ws.Range("H189").Value = 0
$ws.Range("I189").Value = "H"
$ws.Range("J189").Value = 2.25
$ws.Range("L189").Value = 3.5
$ws.Range("M189").Value = 2.625
$ws.Range("N189").Value = 2.875
$ws.Range("O189").Value = 3
$ws.Range("P189").Value = 0
$ws.Range("Q189").Value = 1.8
$ws.Range("R189").Value = 2.05
$ws.Range("S189").Value = 2.25
$ws.Range("T189").Value = 1.975
$ws.Range("U189").Value = 1.875
$ws.Range("V189").Value = 1.625
$ws.Range("X189").Value = -1
$ws.Range("Y189").Value = 0.8
$ws.Range("Z189").Value = -1
$ws.Range("AA189").Value = -1
$ws.Range("AB189").Value = 0.875

# Row 193
$ws.Range("M193").Value = 3.1
$ws.Range("N193").Value = 3.25
$ws.Range("O193").Value = 2.05
$ws.Range("Q193").Value = 2
$ws.Range("R193").Value = 1.85
$ws.Range("S193").Value = 2.5
$ws.Range("T193").Value = 1.825
$ws.Range("U193").Value = 2.025

# Row 194
$ws.Range("M194").Value = 2.3
$ws.Range("O194").Value = 2.8
$ws.Range("P194").Value = -0.25
$ws.Range("Q194").Value = 2.1
$ws.Range("R194").Value = 1.775
$ws.Range("T194").Value = 1.825
$ws.Range("U194").Value = 2.025

# Row 195
$ws.Range("T195").Value = 1.85
$ws.Range("U195").Value = 2

# Row 196
$ws.Range("M196").Value = 2.375
$ws.Range("O196").Value = 2.75
$ws.Range("P196").Value = 0
$ws.Range("Q196").Value = 1.775
$ws.Range("R196").Value = 2.1
